$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 20) with a fresh reading pulled from
# Adafruit IO, mirroring the existing rows where every value --
# including the numeric-looking "Value" column -- is stored as text.
$row = 20

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C ("Value") holds plain numbers as text throughout the sheet.
# Briefly force text formatting so Excel doesn't auto-convert "25" to a
# numeric value, then restore the default "Normal" style so the cell
# doesn't end up with a lingering custom number format.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).NumberFormat = "General"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
